$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices and 1h volume deltas) per latest scrape
# Row 2: update D, E
$ws.Range("D2").Value = '''64.571.73'
$ws.Range("E2").Value = '  +0.13%  '

# Row 3: update D, E
$ws.Range("D3").Value = '''3.427.75'
$ws.Range("E3").Value = '  -1.64%  '

# Row 4: update E
$ws.Range("E4").Value = '  -0.08%  '

# Row 5: update D, E
$ws.Range("D5").Value = '''572.65'
$ws.Range("E5").Value = '  -1.84%  '

# Row 6: update D, E
$ws.Range("D6").Value = '''159.96'
$ws.Range("E6").Value = '  +0.67%  '

# Row 7: update D, E
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.05%  '

# Row 8: update D, E
$ws.Range("D8").Value = '''3.428.37'
$ws.Range("E8").Value = '  -1.94%  '

# Row 9: update D
$ws.Range("D9").Value = '''0.576'

# Row 10: update E
$ws.Range("E10").Value = '  -4.65%  '

# Row 11: update E
$ws.Range("E11").Value = '  -0.38%  '

# Row 12: update D, E
$ws.Range("D12").Value = '''0.436'
$ws.Range("E12").Value = '  -1.64%  '

# Row 13: update D, E
$ws.Range("D13").Value = '''4.012.28'
$ws.Range("E13").Value = '  -1.64%  '

# Row 14: update E
$ws.Range("E14").Value = '  -2.38%  '

# Row 15: update D, E
$ws.Range("D15").Value = '''0.0000192'
$ws.Range("E15").Value = '  +1.59%  '

# Row 16: update D, E
$ws.Range("D16").Value = '''27.97'
$ws.Range("E16").Value = '  +0.60%  '

# Row 17: update D, E
$ws.Range("D17").Value = '''64.571.81'
$ws.Range("E17").Value = '  -0.05%  '

# Row 18: update D, E
$ws.Range("D18").Value = '''3.437.69'
$ws.Range("E18").Value = '  -1.45%  '

# Row 19: update D, E
$ws.Range("D19").Value = '''6.33'
$ws.Range("E19").Value = '  -2.15%  '

# Row 20: update D, E
$ws.Range("D20").Value = '''14.16'
$ws.Range("E20").Value = '  -2.09%  '

# Row 21: update D, E
$ws.Range("D21").Value = '''384.56'
$ws.Range("E21").Value = '  -3.82%  '

# Row 22: update D, E
$ws.Range("D22").Value = '''8.14'
$ws.Range("E22").Value = '  -5.04%  '

# Row 23: update B, C, D, E
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = '''0.541'
$ws.Range("E23").Value = '  -1.21%  '

# Row 24: update B, C, D, E
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''72.70'
$ws.Range("E24").Value = '  +0.77%  '

# Row 25: update E
$ws.Range("E25").Value = '  -0.47%  '

# Row 26: update D, E
$ws.Range("D26").Value = '''0.0000122'
$ws.Range("E26").Value = '  +9.68%  '

# Row 27: update D, E
$ws.Range("D27").Value = '''9.53'
$ws.Range("E27").Value = '  +0.74%  '

# Row 28: update E
$ws.Range("E28").Value = '  -1.86%  '

# Row 29: update E
$ws.Range("E29").Value = '  +0.24%  '

# Row 30: update D, E
$ws.Range("D30").Value = '''6.18'
$ws.Range("E30").Value = '  +3.55%  '

# Row 31: update E
$ws.Range("E31").Value = '  -0.23%  '

# Row 32: update E
$ws.Range("E32").Value = '  -1.31%  '

# Row 33: update B, C, D, E
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''23.52'
$ws.Range("E33").Value = '  -1.70%  '

# Row 34: update B, C, D, E
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '''6.50'
$ws.Range("E34").Value = '  -3.49%  '

# Row 35: update E
$ws.Range("E35").Value = '  +0.12%  '

# Row 36: update D, E
$ws.Range("D36").Value = '''7.07'
$ws.Range("E36").Value = '  +1.40%  '

# Row 37: update D, E
$ws.Range("D37").Value = '''162.08'
$ws.Range("E37").Value = '  +2.12%  '

# Row 38: update D, E
$ws.Range("D38").Value = '''1.49'
$ws.Range("E38").Value = '  -1.66%  '

# Row 39: update D, E
$ws.Range("D39").Value = '''1.91'
$ws.Range("E39").Value = '  +0.80%  '

# Row 40: update D, E
$ws.Range("D40").Value = '''3.013.77'
$ws.Range("E40").Value = '  +4.15%  '

# Row 41: update D, E
$ws.Range("D41").Value = '''0.0762'
$ws.Range("E41").Value = '  -3.74%  '

# Row 42: update D, E
$ws.Range("D42").Value = '''27.08'
$ws.Range("E42").Value = '  -5.77%  '

# Row 43: update D, E
$ws.Range("D43").Value = '''4.51'
$ws.Range("E43").Value = '  +1.48%  '

# Row 44: update B, C, D, E
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0316'
$ws.Range("E44").Value = '  -2.61%  '

# Row 45: update B, C, D, E
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''42.67'
$ws.Range("E45").Value = '  +0.93%  '

# Row 46: update D, E
$ws.Range("D46").Value = '''0.767'
$ws.Range("E46").Value = '  -2.45%  '

# Row 47: update D, E
$ws.Range("D47").Value = '''24.49'
$ws.Range("E47").Value = '  +6.97%  '

# Row 48: update E
$ws.Range("E48").Value = '  -3.44%  '

# Row 49: update D, E
$ws.Range("D49").Value = '''0.868'
$ws.Range("E49").Value = '  +3.01%  '

# Row 50: update D, E
$ws.Range("D50").Value = '''6.60'
$ws.Range("E50").Value = '  +2.57%  '

# Row 51: update B, C, D, E
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '''2.16'
$ws.Range("E51").Value = '  +2.09%  '
